$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C5) from 45183 to 45184 (one day later)
$ws.Range("C2:C5").Value = 45184
